# "Crit bonus changed slightly" — +1 damage when phys, special effect otherwise
# On the "Ranged Weapons P" sheet (sheet2), the crit-bonus-damage input (U5)
# drops from 2 to 0, and the W11 roll value bumps from 4 to 5. All the other
# changed cells in the diff are cached formula results that recalc on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ranged Weapons P")

$ws.Range("U5").Value = 0
$ws.Range("W11").Value = 5

# Restore the author's final active-cell selection on that sheet.
$ws.Activate()
$ws.Range("T4").Select()
